$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 72) for the "Cycling" cell ontology entry.
$ws.Range("A72").Value = "Cycling"
$ws.Range("B72").Value = "cycl"
$ws.Range("C72").Value = "cycling"
$ws.Range("D72").Value = "cycling"
$ws.Range("E72").Value = "cycling"
$ws.Range("F72").Value = "cycling"
